$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.794.38"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.084.81"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'234.02"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "'58.75"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "'0.0790"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "2.394.14"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "'14.74"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "'21.23"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "'0.768"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "2.083.59"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "37.704.55"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "'71.42"
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'228.74"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'170.40"
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("E27").Value = "  +8.15%  "
$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'19.55"
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").Value = "'4.70"
$ws.Range("E32").Value = "  +3.85%  "
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").Value = "'0.0630"
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "'3.46"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'5.41"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").Value = "'0.0983"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'98.77"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "1.459.41"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("E47").Value = "  +6.88%  "
$ws.Range("E48").Value = "  +4.55%  "
$ws.Range("D49").Value = "'7.44"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").Value = "'47.23"
$ws.Range("E51").Value = "  +3.96%  "
